$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Remaining" counts for the Medium section (BackTracking and Math rows)
$ws.Range("D16").Value = 1
$ws.Range("D20").Value = 3

# Move the active selection to D17 (Trees row, Medium section)
$ws.Range("D17").Select()
